# Apply the edit described in the diff:
#  - G2:G4 label changes from "mom1" to "moment"
#  - Two new columns (H: "mom", I: "pol") added with header + three rows of 1s
#  - Selection moves to I7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells H1 ("mom") and I1 ("pol"), styled like the other headers ---
# (introduced before the "moment" rewrite below so new shared strings land in
#  the same order as the target workbook: mom, moment, pol)
$ws.Range("H1").Value = "mom"
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").Font.Bold = $true

# --- Update the existing "mom1" -> "moment" label in column G (rows 2-4) ---
$ws.Range("G2").Value = "moment"
$ws.Range("G3").Value = "moment"
$ws.Range("G4").Value = "moment"

$ws.Range("I1").Value = "pol"
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").Font.Bold = $true

# --- New data cells H2:H4 and I2:I4, all set to 1, styled like the other data cells ---
$ws.Range("H2").Value = 1
$ws.Range("H2").HorizontalAlignment = -4108
$ws.Range("I2").Value = 1
$ws.Range("I2").HorizontalAlignment = -4108

$ws.Range("H3").Value = 1
$ws.Range("H3").HorizontalAlignment = -4108
$ws.Range("I3").Value = 1
$ws.Range("I3").HorizontalAlignment = -4108

$ws.Range("H4").Value = 1
$ws.Range("H4").HorizontalAlignment = -4108
$ws.Range("I4").Value = 1
$ws.Range("I4").HorizontalAlignment = -4108

# --- Move the active selection to I7 (matches the saved sheetView state) ---
$ws.Range("I7").Select()
